$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D4/E4 row: the "removed everything" note is replaced with "None "
$ws.Range("D4").Value = "None "
$ws.Range("E4").Value = "None"

# Add row 5: V1_2 run - simple column drop
$ws.Range("A5").Value = "LinearRegresion"
$ws.Range("B5").Value = "V1_2"
$ws.Range("C5").Value = 0.6107
$ws.Range("D5").Value = "Xóa các cột bị thiếu"
$ws.Range("E5").Value = "None"

# Add row 6: V2 run - drop >50% missing, median-fill the rest
$ws.Range("A6").Value = "LinearRegresion"
$ws.Range("B6").Value = "V2"
$ws.Range("C6").Value = 0.3586
$ws.Range("D6").Value = "Xóa các cột có dữ liệu thiếu trên 50%,điền trung vị cho các cột có dữ liệu dưới 50% bị thiếu"
$ws.Range("E6").Value = "None"

# Column D needs to widen to fit the new, longer text (bestFit autosize)
$ws.Columns("D:D").ColumnWidth = 29.6

# Leave the same selection Excel would land on after typing into E6
$ws.Range("E6").Select() | Out-Null
